$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply AutoFilter on the "Plataform" column (3rd column of Table2) to show only "Mobile" rows ---
# colId="2" in the table XML corresponds to the 3rd field of the table (Plataform),
# which is field index 3 when calling Range.AutoFilter. Operator 7 = xlFilterValues,
# producing a <filters><filter val="Mobile"/></filters> criteria (discrete value filter)
# and hiding every row whose Plataform is not "Mobile".
$lo = $ws.ListObjects.Item(1)
$lo.Range.AutoFilter(3, "Mobile", 7)

# --- Update the "Porcentagem" (column I) values for completed/updated backlog items ---
$ws.Range("I6").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("I14").Value = 100
$ws.Range("I21").Value = 100
$ws.Range("I34").Value = 100
$ws.Range("I44").Value = 100
$ws.Range("I45").Value = 100
$ws.Range("I46").Value = 100
$ws.Range("I47").Value = 100
$ws.Range("I48").Value = 80
$ws.Range("I49").Value = 100

# --- Update the active selection / scroll position on the sheet view ---
$ws.Range("C42").Select()
